$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Port Pins"

# Header row (unchanged values, ensure they are set)
$ws.Range("A1").Value = "Port"
$ws.Range("B1").Value = "Pin"
$ws.Range("C1").Value = "Sub-module"

# Row 2 - Encoder
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "Encoder"
$ws.Range("D2").Value = "GPIO Input (Interrupt)"

# Row 3 - Encoder
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "Encoder"
$ws.Range("D3").Value = "GPIO Input (Interrupt)"

# Row 4 - Barcode Scanner
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "Barcode Scanner"
$ws.Range("D4").Value = "ADC Input"

# Row 6 - Ultrasonic Trigger
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = "Ultrasonic"
$ws.Range("D6").Value = "Trigger"

# Row 7 - Ultrasonic Echo
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = "Ultrasonic"
$ws.Range("D7").Value = "Echo (Interrupt)"

# Row 8 - Ultrasonic LED
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = "Ultrasonic"
$ws.Range("D8").Value = "LED"

# Column widths matching bestFit sizing observed in target
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 18.5

# Selection matches target activeCell
$ws.Range("C17").Select()

$wb.Save()
